## "Big Updates to Query 5"
## Adds a second Unique_Skills / Total_Count summary table (rows 33-41) to the
## "Insights" sheet, mirroring the bottom portion (rows 23-29 -> J2:K8) of the
## helper table already computed on the "Raw Data" sheet, and updates the
## selection/view state on both sheets.

$wb = $excel.ActiveWorkbook
$wsInsights = $wb.Worksheets.Item("Insights")
$wsRaw = $wb.Worksheets.Item("Raw Data")

# ---------------------------------------------------------------------------
# "Insights" sheet: append a second summary table below the existing one.
# ---------------------------------------------------------------------------

# Touch A1 (a no-op format re-assignment) so the sheet's used-range / dimension
# keeps including the blank formatted row 1 above the first table once the
# sheet is re-saved.
$wsInsights.Range("A1").Font.Bold = $false

# Row 33 is a blank spacer row (matches the formatting of row 1 above the
# first table).
$wsInsights.Rows.Item(33).RowHeight = 15.75

# Row 34 holds the new table's header ("Unique_Skills" / "Total_Count"),
# formatted the same way as the header of the first table (row 2, columns A:B).
$wsInsights.Range("A2:B2").Copy() | Out-Null
$wsInsights.Range("A34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsInsights.Rows.Item(34).RowHeight = 15.75
$wsInsights.Range("A34").Value = "Unique_Skills"
$wsInsights.Range("B34").Value = "Total_Count"

# Rows 35-41 hold the new table's data rows, formatted the same way as the
# existing data rows (e.g. row 4, columns A:B).
$wsInsights.Range("A4:B4").Copy() | Out-Null
$wsInsights.Range("A35:B41").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$newSkillCounts = @(
  @("pandas",    3),
  @("excel",     3),
  @("snowflake", 3),
  @("r",         4),
  @("tableau",   6),
  @("python",    7),
  @("sql",       8)
)

$targetRow = 35
foreach ($entry in $newSkillCounts) {
    $wsInsights.Cells.Item($targetRow, 1).Value = $entry[0]
    $wsInsights.Cells.Item($targetRow, 2).Value = $entry[1]
    $targetRow++
}

# ---------------------------------------------------------------------------
# View / selection state.
# ---------------------------------------------------------------------------

# On "Raw Data", the last selection moves to G41.
$wsRaw.Activate()
$wsRaw.Range("G41").Select() | Out-Null

# "Insights" remains the active (displayed) sheet, selection reset to A1.
$wsInsights.Activate()
$wsInsights.Range("A1").Select() | Out-Null
